$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: change fill style from "yellow" (style index 2) to "green" (style index 1)
$ws.Range("A5:B5").Interior.Color = 5296274

# Row 6: apply the same "green" fill style (previously unstyled)
$ws.Range("A6:B6").Interior.Color = 5296274

# New row 13: Chris  | Separation - calculateDistantance test
$ws.Range("A13").Value = "Chris "
$ws.Range("B13").Value = "Separation - calculateDistantance test"

# Move the active selection to A12 (matches the saved selection state of the edited file)
$ws.Range("A12").Select()
